# Weekly update: shift existing rows down by one and insert the new
# week's Jengibre price entry at the top of the data block (row 7).
# Row 66 is a brand-new row, so its "constant" columns (A,B,C,E,F,G,H,I,N,O,Q,R)
# are populated explicitly; rows 7-65 already contain those constant values
# and only need their per-week figures (D,J,K,L,M,P) rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row=7; D=45282; J=50; K=26000; L=26000; M=26000; P=2000 },
    @{ Row=8; D=44467; J=100; K=13000; L=14000; M=13500; P=1038 },
    @{ Row=9; D=45135; J=30; K=18000; L=18000; M=18000; P=1385 },
    @{ Row=10; D=44474; J=40; K=13000; L=14000; M=13500; P=1038 },
    @{ Row=11; D=45126; J=30; K=16000; L=16000; M=16000; P=1231 },
    @{ Row=12; D=44350; J=40; K=23000; L=25000; M=24000; P=1846 },
    @{ Row=13; D=45160; J=40; K=17500; L=18000; M=17750; P=1365 },
    @{ Row=14; D=44782; J=40; K=13000; L=14000; M=13500; P=1038 },
    @{ Row=15; D=45146; J=40; K=17000; L=18000; M=17500; P=1346 },
    @{ Row=16; D=45062; J=30; K=16000; L=17000; M=16333; P=1256 },
    @{ Row=17; D=44503; J=35; K=15000; L=16000; M=15429; P=1187 },
    @{ Row=18; D=45055; J=60; K=18000; L=18000; M=18000; P=1385 },
    @{ Row=19; D=45106; J=50; K=15000; L=16000; M=15600; P=1200 },
    @{ Row=20; D=44453; J=50; K=14000; L=15000; M=14600; P=1123 },
    @{ Row=21; D=44708; J=50; K=13000; L=14000; M=13600; P=1046 },
    @{ Row=22; D=44610; J=50; K=17000; L=18000; M=17400; P=1338 },
    @{ Row=23; D=44433; J=100; K=13000; L=14000; M=13500; P=1038 },
    @{ Row=24; D=45127; J=35; K=17000; L=18000; M=17429; P=1341 },
    @{ Row=25; D=44755; J=40; K=14000; L=15000; M=14500; P=1115 },
    @{ Row=26; D=44509; J=100; K=15000; L=16000; M=15500; P=1192 },
    @{ Row=27; D=44775; J=20; K=12000; L=13000; M=12500; P=962 },
    @{ Row=28; D=44819; J=50; K=13000; L=14000; M=13400; P=1031 },
    @{ Row=29; D=45083; J=30; K=18000; L=18000; M=18000; P=1385 },
    @{ Row=30; D=44355; J=60; K=18000; L=20000; M=19000; P=1462 },
    @{ Row=31; D=44719; J=50; K=13000; L=14000; M=13400; P=1031 },
    @{ Row=32; D=45037; J=50; K=16000; L=17000; M=16400; P=1262 },
    @{ Row=33; D=45125; J=30; K=16000; L=16000; M=16000; P=1231 },
    @{ Row=34; D=44523; J=40; K=15000; L=16000; M=15500; P=1192 },
    @{ Row=35; D=44691; J=100; K=12000; L=13000; M=12500; P=962 },
    @{ Row=36; D=44883; J=60; K=14000; L=15000; M=14500; P=1115 },
    @{ Row=37; D=44813; J=50; K=13000; L=14000; M=13400; P=1031 },
    @{ Row=38; D=44664; J=50; K=11000; L=12000; M=11600; P=892 },
    @{ Row=39; D=45020; J=40; K=15000; L=16000; M=15500; P=1192 },
    @{ Row=40; D=44377; J=40; K=14000; L=15000; M=14500; P=1115 },
    @{ Row=41; D=45013; J=220; K=15000; L=16000; M=15455; P=1189 },
    @{ Row=42; D=45044; J=50; K=20000; L=20000; M=20000; P=1538 },
    @{ Row=43; D=44308; J=50; K=26000; L=27000; M=26400; P=2031 },
    @{ Row=44; D=45167; J=20; K=14000; L=14000; M=14000; P=1077 },
    @{ Row=45; D=45090; J=50; K=15000; L=16000; M=15600; P=1200 },
    @{ Row=46; D=44769; J=50; K=14000; L=15000; M=14600; P=1123 },
    @{ Row=47; D=44313; J=50; K=25000; L=26000; M=25600; P=1969 },
    @{ Row=48; D=44510; J=40; K=15000; L=16000; M=15500; P=1192 },
    @{ Row=49; D=45154; J=50; K=18000; L=18000; M=18000; P=1385 },
    @{ Row=50; D=44327; J=50; K=24000; L=25000; M=24400; P=1877 },
    @{ Row=51; D=44425; J=60; K=14000; L=15000; M=14500; P=1115 },
    @{ Row=52; D=44741; J=50; K=14000; L=15000; M=14400; P=1108 },
    @{ Row=53; D=44362; J=40; K=15000; L=16000; M=15500; P=1192 },
    @{ Row=54; D=44462; J=60; K=14000; L=15000; M=14500; P=1115 },
    @{ Row=55; D=44316; J=50; K=27000; L=28000; M=27400; P=2108 },
    @{ Row=56; D=44383; J=50; K=15000; L=16000; M=15400; P=1185 },
    @{ Row=57; D=44334; J=50; K=26000; L=28000; M=27200; P=2092 },
    @{ Row=58; D=44978; J=40; K=13000; L=14000; M=13500; P=1038 },
    @{ Row=59; D=44761; J=25; K=14000; L=15000; M=14400; P=1108 },
    @{ Row=60; D=44777; J=25; K=13000; L=14000; M=13600; P=1046 },
    @{ Row=61; D=44810; J=50; K=11000; L=12000; M=11600; P=892 },
    @{ Row=62; D=45142; J=30; K=18000; L=18000; M=18000; P=1385 },
    @{ Row=63; D=44488; J=40; K=16000; L=17000; M=16500; P=1269 },
    @{ Row=64; D=44705; J=50; K=10000; L=11000; M=10400; P=800 },
    @{ Row=65; D=45118; J=30; K=15000; L=15000; M=15000; P=1154 },
    @{ Row=66; D=45155; J=25; K=15000; L=15000; M=15000; P=1154 }

)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value = $item.D    # D: Fecha
    $ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat
    $ws.Cells.Item($r, 10).Value = $item.J   # J: Volumen
    $ws.Cells.Item($r, 11).Value = $item.K   # K: Precio minimo
    $ws.Cells.Item($r, 12).Value = $item.L   # L: Precio maximo
    $ws.Cells.Item($r, 13).Value = $item.M   # M: Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $item.P   # P: Precio $/Kg
}

# Row 66 is entirely new, so fill in the columns that stay constant across
# every data row (copied verbatim from an existing row, e.g. row 65).
# NOTE: use .Value2 (not .Value) when *reading* a cell back in this engine.
$ws.Cells.Item(66, 1).Value = $ws.Cells.Item(65, 1).Value2    # A: Mercado ID
$ws.Cells.Item(66, 2).Value = $ws.Cells.Item(65, 2).Value2    # B: Mercado
$ws.Cells.Item(66, 3).Value = $ws.Cells.Item(65, 3).Value2    # C: Region
$ws.Cells.Item(66, 5).Value = $ws.Cells.Item(65, 5).Value2    # E: Codreg
$ws.Cells.Item(66, 6).Value = $ws.Cells.Item(65, 6).Value2    # F: Categoria ID
$ws.Cells.Item(66, 7).Value = $ws.Cells.Item(65, 7).Value2    # G: Categoria
$ws.Cells.Item(66, 8).Value = $ws.Cells.Item(65, 8).Value2    # H: Variedad
$ws.Cells.Item(66, 9).Value = $ws.Cells.Item(65, 9).Value2    # I: Calidad
$ws.Cells.Item(66, 14).Value = $ws.Cells.Item(65, 14).Value2  # N: Unidad de comercializacion
$ws.Cells.Item(66, 15).Value = $ws.Cells.Item(65, 15).Value2  # O: Origen
$ws.Cells.Item(66, 17).Value = $ws.Cells.Item(65, 17).Value2  # Q: Kg o Unidades
$ws.Cells.Item(66, 18).Value = $ws.Cells.Item(65, 18).Value2  # R: Clasificacion
